# The post at row 528 ("「友達はハム太郎のことが好き」") was removed from the
# workbook. Deleting its entire row shifts every subsequent row up by one,
# which matches the target diff (row 529 -> 528, 530 -> 529, ..., 661 -> 660)
# and also updates the sheet's used-range dimension from C661 to C660
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(528).Delete()
